$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$values = @(
    @(9, 9),
    @(8, 8),
    @(12, 12),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(5, 5),
    @(4, 4),
    @(5, 5),
    @(3, 3),
    @(2, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
